# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.741.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.369.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.366.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.796.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.647.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.369.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.40%  "
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.485.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "520.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.46%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0907"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.45%  "
